$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: the paragraph ending in "... Zie screendump hieronder."
# gets its trailing run " hieronder." split into three runs:
#   " "  /  "hieronder"  /  "."
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Zie screendump hieronder.*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start
$full = $target.Range.Text
$needle = " hieronder."
$idx = $full.IndexOf($needle)
$absStart = $pStart + $idx
$absEnd = $absStart + $needle.Length
$sub = $d.Range($absStart, $absEnd)
$sub.Text = " "
$sub.Collapse(0)
$sub.InsertAfter("hieronder")
$sub.Collapse(0)
$sub.InsertAfter(".")

# ------------------------------------------------------------------
# Step 2: insert a brand new "List Paragraph" styled paragraph right
# after it, describing the SpectreYggdra / ASIX adapter choice.
# ------------------------------------------------------------------
$target = $null
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Zie screendump hieronder.*") {
        $target = $p
        $targetIndex = $i
        break
    }
}

$r = $target.Range
$r.End = $r.End - 1
$r.Collapse(0)
$newPara = $r.InsertParagraphAfter()

$newParaObj = $d.Paragraphs($targetIndex + 1)
$newParaObj.Range.Style = "List Paragraph"

$ir = $d.Paragraphs($targetIndex + 1).Range
$ir.End = $ir.End - 1
$ir.Collapse(0)

$ir.InsertAfter("Voor")
$ir.Collapse(0)
$ir.InsertAfter(" de ")
$ir.Collapse(0)
$ir.InsertAfter("SpectreYggdra")
$ir.Collapse(0)
$ir.InsertAfter(" (pc Frank van Bokhoven) ")
$ir.Collapse(0)
$ir.InsertAfter("moet")
$ir.Collapse(0)
$ir.InsertAfter(" ")
$ir.Collapse(0)
$ir.InsertAfter("gekozen")
$ir.Collapse(0)
$ir.InsertAfter(" ")
$ir.Collapse(0)
$ir.InsertAfter("worden")
$ir.Collapse(0)
$ir.InsertAfter(" ")
$ir.Collapse(0)
$ir.InsertAfter("voor")
$ir.Collapse(0)
$ir.InsertAfter(" de ASIX AX88772A USB 2.0 to Fast Ethernet Adapter.")
$ir.Collapse(0)
$ir.InsertAfter(" (")
$ir.Collapse(0)
$ir.InsertAfter("zie")
$ir.Collapse(0)
$ir.InsertAfter(" ")
$ir.Collapse(0)
$ir.InsertAfter("screendump")
$ir.Collapse(0)
$ir.InsertAfter("). ")
$ir.Collapse(0)
$ir.InsertAfter("Dit")
$ir.Collapse(0)
$ir.InsertAfter(" is de adapter die ")
$ir.Collapse(0)
$ir.InsertAfter("verbonden")
$ir.Collapse(0)
$ir.InsertAfter(" is met de router ")
$ir.Collapse(0)
$ir.InsertAfter("waaraan")
$ir.Collapse(0)
$ir.InsertAfter(" ")
$ir.Collapse(0)
$ir.InsertAfter("weer")
$ir.Collapse(0)
$ir.InsertAfter(" de CT’s ")
$ir.Collapse(0)
$ir.InsertAfter("verbonden")
$ir.Collapse(0)
$ir.InsertAfter(" zijn.")
$ir.Collapse(0)

Write-Host "Edit complete."
